$wb = $excel.ActiveWorkbook

# ======================================================================
# Sheet: GradeOneVideoList  (xl/worksheets/sheet4.xml)
# ======================================================================
$ws = $wb.Worksheets.Item("GradeOneVideoList")

# ---- Rows 2-10: replace existing data ----
$ws.Range("A2").Value = "Bible"
$ws.Range("B2").Value = "Lesson 1"
$ws.Range("C2").Value = "Reading 1"
$ws.Range("D2").Value = "Reading 1"
$ws.Range("D2").Style = "Normal"

$ws.Range("A3").Value = "Classroom Routines"
$ws.Range("B3").Value = "Lesson 1"
$ws.Range("C3").Value = "Bible 1"
$ws.Range("D3").Value = "Bible 1"
$ws.Range("D3").Style = "Normal"

$ws.Range("A4").Value = "Seatwork Explanation (Cursive)"
$ws.Range("B4").Value = "Lesson 1"
$ws.Range("C4").Value = "Spelling 1"
$ws.Range("D4").Value = "Spelling 1"
$ws.Range("D4").Style = "Normal"

$ws.Range("A5").Value = "Phonics/Language"
$ws.Range("B5").Value = "Lesson 1"
$ws.Range("C5").Value = "Phonics 1"
$ws.Range("D5").Value = "Phonics 1"
$ws.Range("D5").Style = "Normal"

$ws.Range("A6").Value = "Cursive Writing"
$ws.Range("B6").Value = "Lesson 1"
$ws.Range("C6").Value = "Activities 1"
$ws.Range("D6").Value = "Activities 1"
$ws.Range("D6").Style = "Normal"

$ws.Range("A7").Value = "Spelling/Poetry"
$ws.Range("B7").Value = "Lesson 1"
$ws.Range("C7").Value = "Writing 1"
$ws.Range("D7").Value = "Writing 1"
$ws.Range("D7").Style = "Normal"
$ws.Range("F7").Value = 1

$ws.Range("A8").Value = "Arithmetic"
$ws.Range("B8").Value = "Lesson 1"
$ws.Range("C8").Value = "Seatwork 1"
$ws.Range("D8").Value = "Seatwork 1"
$ws.Range("D8").Style = "Normal"

$ws.Range("A9").Value = "Combination Practice"
$ws.Range("B9").Value = "Lesson 1"
$ws.Range("C9").Value = "Arithmetic 1"
$ws.Range("D9").Value = "Arithmetic 1"
$ws.Range("D9").Style = "Normal"

$ws.Range("A10").Value = "Activity Time"
$ws.Range("B10").Value = "Lesson 1"
$ws.Range("C10").Value = "Classroom Routines 1"
$ws.Range("D10").Value = "Classroom Routines 1"
$ws.Range("D10").Style = "Normal"

# ---- Row 19: new lone label ----
$ws.Range("A19").Value = "Manuscript Subject"

# ---- Row 20: new header row (copy style from row 1's header) ----
$ws.Range("A1:G1").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122)
$ws.Range("A20").Value = "MyLessonsTodaySubjectList"
$ws.Range("B20").Value = "MyLessonsTodayLessonList"
$ws.Range("C20").Value = "VideoLibraryDropdownSubjectList"
$ws.Range("D20").Value = "VideoLibraryDropdownLongDescriptionList"
$ws.Range("E20").Value = "SegmentId"
$ws.Range("F20").Value = "TodayLessonOfVideoLibrary"
$ws.Range("G20").Value = "NextDayLessonOfVideoLibrary"

# ---- Rows 21-29: new data block ----
$ws.Range("A21").Value = "Bible"
$ws.Range("B21").Value = "Lesson 1"
$ws.Range("C21").Value = "Reading 1"
$ws.Range("D21").Value = "Reading 1"
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 2

$ws.Range("A22").Value = "Classroom Routines"
$ws.Range("B22").Value = "Lesson 1"
$ws.Range("C22").Value = "Bible 1"
$ws.Range("D22").Value = "Bible 1"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2

$ws.Range("A23").Value = "Seatwork Explanation (Manuscript)"
$ws.Range("B23").Value = "Lesson 1"
$ws.Range("C23").Value = "Spelling 1"
$ws.Range("D23").Value = "Spelling 1"
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 2

$ws.Range("A24").Value = "AM Giraffes"
$ws.Range("B24").Value = "Lesson 4"
$ws.Range("C24").Value = "Phonics 1"
$ws.Range("D24").Value = "Phonics 1"
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 2

$ws.Range("A25").Value = "AM Monkeys"
$ws.Range("B25").Value = "Lesson 4"
$ws.Range("C25").Value = "Activities 1"
$ws.Range("D25").Value = "Activities 1"
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 2

$ws.Range("A26").Value = "AM Elephants"
$ws.Range("B26").Value = "Lesson 4"
$ws.Range("C26").Value = "Writing 1"
$ws.Range("D26").Value = "Writing 1"
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 2

$ws.Range("A27").Value = "Phonics/Language"
$ws.Range("B27").Value = "Lesson 1"
$ws.Range("C27").Value = "Seatwork 1"
$ws.Range("D27").Value = "Seatwork 1"
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 2

$ws.Range("A28").Value = "Manuscript Writing"
$ws.Range("B28").Value = "Lesson 1"
$ws.Range("C28").Value = "Arithmetic 1"
$ws.Range("D28").Value = "Arithmetic 1"
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 2

$ws.Range("A29").Value = "Spelling/Poetry"
$ws.Range("B29").Value = "Lesson 1"
$ws.Range("C29").Value = "Classroom Routines 1"
$ws.Range("D29").Value = "Classroom Routines 1"
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 2

# ---- Rows 30-32: only A & B columns ----
$ws.Range("A30").Value = "Combination Practice"
$ws.Range("B30").Value = "Lesson 1"

$ws.Range("A31").Value = "Arithmetic"
$ws.Range("B31").Value = "Lesson 1"

$ws.Range("A32").Value = "Activity Time"
$ws.Range("B32").Value = "Lesson 1"

# ---- New bestFit columns C & D ----
$ws.Columns.Item(3).ColumnWidth = 26.8
$ws.Columns.Item(4).ColumnWidth = 33.6

# ======================================================================
# Sheet: StudentCredentials  (xl/worksheets/sheet2.xml) - selection change
# ======================================================================
$ws2 = $wb.Worksheets.Item("StudentCredentials")

# ======================================================================
# Sheet: GradeNineVideoList (xl/worksheets/sheet5.xml) - was the active tab
# ======================================================================
$ws5 = $wb.Worksheets.Item("GradeNineVideoList")

# Move active tab from GradeNineVideoList to GradeOneVideoList, and set
# selections for the touched sheets (do this last so later value writes
# above don't move the selection away from the intended cell).
$ws2.Activate()
$ws2.Range("A3").Select()

$ws.Activate()
$ws.Range("C21").Select()
